# Weekly refresh of the "Chirimoya" price sheet: the per-market price
# records (date, quality, volume, min/max/avg price, unit, $/kg, kg/unit)
# were reshuffled across rows 2-14 while keeping the static market/product
# metadata columns (A,B,C,E,F,G,H,I,J,K,R) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 5's data pattern)
$ws.Cells.Item(2,4).Value = 44482
$ws.Cells.Item(2,12).Value = "Primera"
$ws.Cells.Item(2,13).Value = 160
$ws.Cells.Item(2,14).Value = 25000
$ws.Cells.Item(2,15).Value = 26000
$ws.Cells.Item(2,16).Value = 25500
$ws.Cells.Item(2,19).Value = 2125

# Row 3 (was row 2's data pattern)
$ws.Cells.Item(3,4).Value = 44475
$ws.Cells.Item(3,12).Value = "Especial"
$ws.Cells.Item(3,14).Value = 32000
$ws.Cells.Item(3,15).Value = 33000
$ws.Cells.Item(3,16).Value = 32500
$ws.Cells.Item(3,17).Value = "$/caja 12 kilos"
$ws.Cells.Item(3,19).Value = 2708
$ws.Cells.Item(3,20).Value = 12

# Row 4 (was row 14's data pattern)
$ws.Cells.Item(4,4).Value = 44489
$ws.Cells.Item(4,13).Value = 200
$ws.Cells.Item(4,14).Value = 24000
$ws.Cells.Item(4,15).Value = 25000
$ws.Cells.Item(4,16).Value = 24500
$ws.Cells.Item(4,19).Value = 2042

# Row 5 (was row 12's data pattern)
$ws.Cells.Item(5,4).Value = 44524
$ws.Cells.Item(5,13).Value = 200
$ws.Cells.Item(5,14).Value = 23000
$ws.Cells.Item(5,15).Value = 24000
$ws.Cells.Item(5,16).Value = 23500
$ws.Cells.Item(5,19).Value = 1958

# Row 6 (was row 7's data pattern)
$ws.Cells.Item(6,4).Value = 44783
$ws.Cells.Item(6,12).Value = "Tercera"
$ws.Cells.Item(6,13).Value = 100
$ws.Cells.Item(6,14).Value = 27000
$ws.Cells.Item(6,15).Value = 28000
$ws.Cells.Item(6,16).Value = 27500
$ws.Cells.Item(6,17).Value = "$/caja 12 kilos"
$ws.Cells.Item(6,19).Value = 2292
$ws.Cells.Item(6,20).Value = 12

# Row 7 (was row 3's data pattern)
$ws.Cells.Item(7,4).Value = 44167
$ws.Cells.Item(7,12).Value = "Segunda"
$ws.Cells.Item(7,13).Value = 200
$ws.Cells.Item(7,14).Value = 18000
$ws.Cells.Item(7,15).Value = 19000
$ws.Cells.Item(7,16).Value = 18500
$ws.Cells.Item(7,17).Value = "$/caja 13 kilos"
$ws.Cells.Item(7,19).Value = 1423
$ws.Cells.Item(7,20).Value = 13

# Row 8 (was row 6's data pattern)
$ws.Cells.Item(8,4).Value = 44468
$ws.Cells.Item(8,12).Value = "Primera"
$ws.Cells.Item(8,14).Value = 29000
$ws.Cells.Item(8,15).Value = 30000
$ws.Cells.Item(8,16).Value = 29500
$ws.Cells.Item(8,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(8,19).Value = 2950
$ws.Cells.Item(8,20).Value = 10

# Row 9 (was row 13's data pattern)
$ws.Cells.Item(9,4).Value = 44441
$ws.Cells.Item(9,12).Value = "Primera"
$ws.Cells.Item(9,13).Value = 100
$ws.Cells.Item(9,17).Value = "$/caja 12 kilos"
$ws.Cells.Item(9,19).Value = 2458
$ws.Cells.Item(9,20).Value = 12

# Row 10 (was row 8's data pattern)
$ws.Cells.Item(10,4).Value = 44160
$ws.Cells.Item(10,12).Value = "Segunda"
$ws.Cells.Item(10,14).Value = 19000
$ws.Cells.Item(10,15).Value = 20000
$ws.Cells.Item(10,16).Value = 19500
$ws.Cells.Item(10,17).Value = "$/caja 13 kilos"
$ws.Cells.Item(10,19).Value = 1500
$ws.Cells.Item(10,20).Value = 13

# Row 11 (was row 9's data pattern)
$ws.Cells.Item(11,4).Value = 44776
$ws.Cells.Item(11,12).Value = "Segunda"
$ws.Cells.Item(11,13).Value = 160
$ws.Cells.Item(11,14).Value = 29000
$ws.Cells.Item(11,15).Value = 30000
$ws.Cells.Item(11,16).Value = 29500
$ws.Cells.Item(11,17).Value = "$/caja 10 kilos"
$ws.Cells.Item(11,19).Value = 2950
$ws.Cells.Item(11,20).Value = 10

# Row 12 (was row 10's data pattern)
$ws.Cells.Item(12,4).Value = 44545
$ws.Cells.Item(12,17).Value = "$/bandeja 12 kilos"

# Row 13 (was row 4's data pattern)
$ws.Cells.Item(13,4).Value = 44811

# Row 14 (was row 11's data pattern)
$ws.Cells.Item(14,4).Value = 44496
$ws.Cells.Item(14,14).Value = 23000
$ws.Cells.Item(14,15).Value = 24000
$ws.Cells.Item(14,16).Value = 23500
$ws.Cells.Item(14,19).Value = 1958

